$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

function Set-TextValue($range, [string]$value) {
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

# Enterprises density (per 1000 people) - row 11
Set-TextValue $ws.Range("B11") "40.76"
Set-TextValue $ws.Range("C11") "4.49"
Set-TextValue $ws.Range("D11") "45.25"

# Employment (% of total) - row 12
Set-TextValue $ws.Range("B12") "21.38"
Set-TextValue $ws.Range("C12") "39.15"
Set-TextValue $ws.Range("D12") "60.53"

# Enterprises (% of total) - row 14
Set-TextValue $ws.Range("B14") "89.62"
Set-TextValue $ws.Range("C14") "9.87"

# Source Type: SME Associations block - row 30
Set-TextValue $ws.Range("B30") "31.65"
Set-TextValue $ws.Range("C30") "1.46"
Set-TextValue $ws.Range("D30") "33.11"
